$p = $ppt.ActivePresentation

# Slide 16 contains a table (graphicFrame, shape 3) whose table style is being
# switched from the deck's custom default style to a different (standard)
# table style id.
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{C9B250B1-A9AB-462E-AD17-D7F7EB72DAB4}")
